$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the dataset row (row 6): kdd_cup_1998 -> sp500 ---
$ws.Range("B6").Value = "sp500"
$ws.Range("C6").Value = "https://www.kaggle.com/camnugent/sandp500"
$ws.Range("D6").Value = "Stock market data can be interesting to analyze and as a further incentive, strong predictive models can have large financial payoff. The amount of financial data on the web is seemingly endless. A large and well structured dataset on a wide array of companies can be hard to come by. Here I provide a dataset with historical stock prices (last 5 years) for all companies currently found on the S&P 500 index. "

$ws.Range("E6").Value = 619404
$ws.Range("F6").Value = 7
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 0

# --- Update the active selection ---
$ws.Range("E7").Select()
